# Fixed bug with non-greedy regular expression matching.
# The underlying text-pattern percentages recomputed by the (fixed) regex
# produced a new set of values for row 2 (A2:AI2) of the data sheet.
# Update the sheet with the corrected values and restore the view state
# (active sheet/selection/zoom) that Excel recorded after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = 0.037792895
$ws.Range("B2").Value  = 0.113378685
$ws.Range("C2").Value  = 0
$ws.Range("D2").Value  = 3.099017385
$ws.Range("E2").Value  = 0.15117158
$ws.Range("F2").Value  = 29.40287226
$ws.Range("G2").Value  = 5.555555556
$ws.Range("H2").Value  = 8.65457294
$ws.Range("I2").Value  = 2.872260015
$ws.Range("J2").Value  = 2.9856387
$ws.Range("K2").Value  = 7.218442933
$ws.Range("L2").Value  = 9.486016629
$ws.Range("M2").Value  = 0.453514739
$ws.Range("N2").Value  = 1.171579743
$ws.Range("O2").Value  = 2.456538171
$ws.Range("P2").Value  = 3.741496599
$ws.Range("Q2").Value  = 0.037792895
$ws.Range("R2").Value  = 0.22675737
$ws.Range("S2").Value  = 1.927437642
$ws.Range("T2").Value  = 0.453514739
$ws.Range("U2").Value  = 9.977324263
$ws.Range("V2").Value  = 1.965230537
$ws.Range("W2").Value  = 0.113378685
$ws.Range("X2").Value  = 3.439153439
$ws.Range("Y2").Value  = 1.549508692
$ws.Range("Z2").Value  = 0
$ws.Range("AA2").Value = 0.15117158
$ws.Range("AB2").Value = 1.700680272
$ws.Range("AC2").Value = 0.22675737
$ws.Range("AD2").Value = 0.037792895
$ws.Range("AE2").Value = 0.113378685
$ws.Range("AF2").Value = 0.037792895
$ws.Range("AG2").Value = 0.113378685
$ws.Range("AH2").Value = 0.453514739
$ws.Range("AI2").Value = 0.07558579

# Restore the sheet view state recorded in the workbook after the edit:
# the data worksheet becomes the active/selected tab, the whole data
# range is selected and the zoom level is changed.
$null = $ws.Activate()
$null = $ws.Range("A1:AI2").Select()
$excel.ActiveWindow.Zoom = 43
